$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.600.94'
$ws.Range("E2").Value = '  +1.09%  '
# Row 3
$ws.Range("D3").Value = '1.985.74'
$ws.Range("E3").Value = '  +4.58%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.35%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '
# Row 6
$ws.Range("E6").Value = '  +0.37%  '
# Row 7
$ws.Range("E7").Value = '  +0.40%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3928'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.17%  '
# Row 9
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07938'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.47%  '
# Row 10
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9936'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.36%  '
# Row 11
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.85'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.70%  '
# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '2.001.45'
$ws.Range("E12").Value = '  +5.51%  '
# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.190'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.55%  '
# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.848'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.70%  '
# Row 15
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07093'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.21%  '
# Row 16
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.95%  '
# Row 17
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.43%  '
# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009942'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.68%  '
# Row 19
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.30%  '
# Row 20
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.28%  '
# Row 21
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '29.629.34'
$ws.Range("E21").Value = '  +1.18%  '
# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.555'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.92%  '
# Row 23
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.96%  '
# Row 24
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.228.89'
$ws.Range("E24").Value = '  +4.78%  '
# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.107'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.56%  '
# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.87%  '
# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.62%  '
# Row 28
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.823'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.49%  '
# Row 29
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '119.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.79%  '
# Row 30
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.902'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.15%  '
# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09419'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.46%  '
# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8991'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.78%  '
# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.226'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.41%  '
# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.334'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.37%  '
# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.194'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.07%  '
# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05816'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.25%  '
# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.177'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.59%  '
# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02101'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.30%  '
# Row 39
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.000003273'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +51.07%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.846'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '
# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5733'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.33%  '
# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1802'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.90%  '
# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.731'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.30%  '
# Row 44
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.765'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.77%  '
# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.24%  '
# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5367'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.14%  '
# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.183'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.84%  '
# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06946'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.48%  '
# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.18%  '
# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.826'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.06%  '
# Row 51
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3043'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.81%  '
